# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# New "Periodo Mora" labels for rows 16-26 (column E), now in ascending order
$periodos = @("1811", "1812", "1901", "1902", "1903", "1904", "1905", "1906", "1907", "1908", "1909")

for ($i = 0; $i -lt $periodos.Length; $i++) {
    $row = 16 + $i
    $ws.Range("E$row").Value = $periodos[$i]
}

# Update "Valor Mora" (column G) for rows 16-26 from 1300000 to 1423500
for ($row = 16; $row -le 26; $row++) {
    $ws.Range("G$row").Value = 1423500
}
